# section_list.xlsx: rebuild the sheet's header/data columns for the new
# import-service schema (course_id, section_id, start, end, classroom_no,
# lesson, limit, day) and drop the old columns (title, credits, dept_name,
# time, instructor_id) plus the second sample data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert a blank column at K. This pushes the existing column K
# (instructor_id) -- and, importantly, its custom width formatting -- one
# slot to the right, to column L. That reproduces the workbook's <cols>
# entry ending up at min=12,max=12 once the unused columns are cleared
# below (Excel keeps per-column formatting even past the used data range).
$ws.Columns("K:K").Insert() | Out-Null

# Step 2: clear the columns that are no longer needed (old limit/day/
# instructor_id data plus the blank column just inserted), leaving only
# columns A-H in use.
$ws.Range("I1:L3").ClearContents() | Out-Null

# Step 3: write the new header row.
$ws.Range("A1").Value = "course_id"
$ws.Range("B1").Value = "section_id"
$ws.Range("C1").Value = "start"
$ws.Range("D1").Value = "end"
$ws.Range("E1").Value = "classroom_no"
$ws.Range("F1").Value = "lesson"
$ws.Range("G1").Value = "limit"
$ws.Range("H1").Value = "day"

# Step 4: write the single remaining data row with the new sample section.
$ws.Range("A2").Value = "CCCC120001"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = "Z2204"
$ws.Range("F2").Style = "Normal"
$ws.Range("F2").Value = 4
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "10"
$ws.Range("H2").Value = 3

# Step 5: the old row 3 (a second full data row) goes away, replaced by an
# otherwise-empty row that keeps only a text-formatted G3 cell (matching
# the trailing style left behind in the authored file). Deleting the row
# and inserting a fresh one picks up formatting from row 2 above, so do
# this after row 2's own formatting has been finalized.
$ws.Rows("3:3").Delete() | Out-Null
$ws.Rows("3:3").Insert() | Out-Null
$ws.Range("G3").NumberFormat = "@"

# Step 6: restore the cursor/selection shown in the saved workbook.
$ws.Range("H5").Select() | Out-Null
